$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.2313123333333333
$ws.Range("H2").Value = 0.693937
$ws.Range("I2").Value = 0.7569517164947553
$ws.Range("J2").Value = 0.7569517164947555
$ws.Range("M2").Value = 0.2313123333333333
$ws.Range("N2").Value = 0.693937
$ws.Range("O2").Value = 0.7569517164947553
$ws.Range("P2").Value = 0.7569517164947555
$ws.Range("Q2").Value = 0.05350539555211112
$ws.Range("R2").Value = 0.481548559969
$ws.Range("S2").Value = 0.5729759011043565
$ws.Range("T2").Value = 0.5729759011043566

# Row 3
$ws.Range("G3").Value = 0.2313123333333333
$ws.Range("H3").Value = 0.693937
$ws.Range("I3").Value = 0.7569517164947553
$ws.Range("J3").Value = 0.7569517164947555
$ws.Range("O3").Value = 0.2385552472206224
$ws.Range("P3").Value = 0.2385552472206224
$ws.Range("Q3").Value = 0.01686236068355556
$ws.Range("R3").Value = 0.151761246152
$ws.Range("S3").Value = 0.1805748038624808
$ws.Range("T3").Value = 0.1805748038624808

# Row 4
$ws.Range("G4").Value = 0.2313123333333333
$ws.Range("H4").Value = 0.693937
$ws.Range("I4").Value = 0.7569517164947553
$ws.Range("J4").Value = 0.7569517164947555
$ws.Range("M4").Value = 0.001373
$ws.Range("N4").Value = 0.004119
$ws.Range("O4").Value = 0.004493036284622232
$ws.Range("P4").Value = 0.004493036284622232
$ws.Range("Q4").Value = 0.0003175918336666667
$ws.Range("R4").Value = 0.002858326503
$ws.Range("S4").Value = 0.003401011527918016
$ws.Range("T4").Value = 0.003401011527918017

# Row 5
$ws.Range("I5").Value = 0.2385552472206224
$ws.Range("J5").Value = 0.2385552472206224
$ws.Range("M5").Value = 0.2313123333333333
$ws.Range("N5").Value = 0.693937
$ws.Range("O5").Value = 0.7569517164947553
$ws.Range("P5").Value = 0.7569517164947555
$ws.Range("Q5").Value = 0.01686236068355556
$ws.Range("R5").Value = 0.151761246152
$ws.Range("S5").Value = 0.1805748038624808
$ws.Range("T5").Value = 0.1805748038624808

# Row 6
$ws.Range("I6").Value = 0.2385552472206224
$ws.Range("J6").Value = 0.2385552472206224
$ws.Range("O6").Value = 0.2385552472206224
$ws.Range("P6").Value = 0.2385552472206224
$ws.Range("S6").Value = 0.05690860597649225
$ws.Range("T6").Value = 0.05690860597649226

# Row 7
$ws.Range("I7").Value = 0.2385552472206224
$ws.Range("J7").Value = 0.2385552472206224
$ws.Range("M7").Value = 0.001373
$ws.Range("N7").Value = 0.004119
$ws.Range("O7").Value = 0.004493036284622232
$ws.Range("P7").Value = 0.004493036284622232
$ws.Range("Q7").Value = 0.0001000898693333333
$ws.Range("R7").Value = 0.0009008088240000001
$ws.Range("S7").Value = 0.001071837381649283
$ws.Range("T7").Value = 0.001071837381649283

# Row 8
$ws.Range("G8").Value = 0.001373
$ws.Range("H8").Value = 0.004119
$ws.Range("I8").Value = 0.004493036284622232
$ws.Range("J8").Value = 0.004493036284622232
$ws.Range("M8").Value = 0.2313123333333333
$ws.Range("N8").Value = 0.693937
$ws.Range("O8").Value = 0.7569517164947553
$ws.Range("P8").Value = 0.7569517164947555
$ws.Range("Q8").Value = 0.0003175918336666667
$ws.Range("R8").Value = 0.002858326503
$ws.Range("S8").Value = 0.003401011527918016
$ws.Range("T8").Value = 0.003401011527918017

# Row 9
$ws.Range("G9").Value = 0.001373
$ws.Range("H9").Value = 0.004119
$ws.Range("I9").Value = 0.004493036284622232
$ws.Range("J9").Value = 0.004493036284622232
$ws.Range("O9").Value = 0.2385552472206224
$ws.Range("P9").Value = 0.2385552472206224
$ws.Range("Q9").Value = 0.0001000898693333333
$ws.Range("R9").Value = 0.0009008088240000001
$ws.Range("S9").Value = 0.001071837381649283
$ws.Range("T9").Value = 0.001071837381649283

# Row 10
$ws.Range("G10").Value = 0.001373
$ws.Range("H10").Value = 0.004119
$ws.Range("I10").Value = 0.004493036284622232
$ws.Range("J10").Value = 0.004493036284622232
$ws.Range("M10").Value = 0.001373
$ws.Range("N10").Value = 0.004119
$ws.Range("O10").Value = 0.004493036284622232
$ws.Range("P10").Value = 0.004493036284622232
$ws.Range("Q10").Value = 0.000001885129
$ws.Range("R10").Value = 0.000016966161
$ws.Range("S10").Value = 0.00002018737505493195
$ws.Range("T10").Value = 0.00002018737505493195
